$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mago")

# Fill in the newly-graded "S4" column (E) for the students that have a
# recorded result for that test.
$ws.Range("E4").Value  = 1
$ws.Range("E5").Value  = 0
$ws.Range("E8").Value  = 1
$ws.Range("E9").Value  = 1
$ws.Range("E10").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("E13").Value = 1
$ws.Range("E14").Value = 1
$ws.Range("E16").Value = 1
$ws.Range("E18").Value = 1
$ws.Range("E19").Value = 1
$ws.Range("E20").Value = 1
$ws.Range("E21").Value = 1

# Make "Mago" the active sheet/tab, with the last-entered cell selected.
$ws.Activate()
$ws.Range("E11").Select()
